# Generate Report for Handoff
# A new handoff round completed for the "9e58fb96-..." file: its "Latest
# Handoff Datetime" is refreshed on both locale status sheets while every
# other recorded value (including the still-pending "a07768e3-..." row)
# is left untouched.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-02-22 17:14:25"
$wsZhCn.Range("D5").Value = "2016-02-22 17:13:45"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-02-22 17:14:36"
$wsDeDe.Range("D5").Value = "2016-02-22 17:13:56"
